$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "Fri Sep 29 11:38:43 EDT 2023"
$ws.Range("B3").Value = "Fri Sep 29 11:38:58 EDT 2023"
$ws.Range("B4").Value = "Fri Sep 29 11:39:13 EDT 2023"
